$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap pairs of rows (content in columns B:AC), keeping column A (sequential index) fixed.
$rowA = $ws.Range("B9:AC9").Value2
$rowB = $ws.Range("B10:AC10").Value2
$ws.Range("B9:AC9").Value2 = $rowB
$ws.Range("B10:AC10").Value2 = $rowA

$rowA = $ws.Range("B17:AC17").Value2
$rowB = $ws.Range("B18:AC18").Value2
$ws.Range("B17:AC17").Value2 = $rowB
$ws.Range("B18:AC18").Value2 = $rowA

$rowA = $ws.Range("B19:AC19").Value2
$rowB = $ws.Range("B20:AC20").Value2
$ws.Range("B19:AC19").Value2 = $rowB
$ws.Range("B20:AC20").Value2 = $rowA

$rowA = $ws.Range("B27:AC27").Value2
$rowB = $ws.Range("B28:AC28").Value2
$ws.Range("B27:AC27").Value2 = $rowB
$ws.Range("B28:AC28").Value2 = $rowA

$rowA = $ws.Range("B38:AC38").Value2
$rowB = $ws.Range("B39:AC39").Value2
$ws.Range("B38:AC38").Value2 = $rowB
$ws.Range("B39:AC39").Value2 = $rowA

$rowA = $ws.Range("B66:AC66").Value2
$rowB = $ws.Range("B67:AC67").Value2
$ws.Range("B66:AC66").Value2 = $rowB
$ws.Range("B67:AC67").Value2 = $rowA

$rowA = $ws.Range("B82:AC82").Value2
$rowB = $ws.Range("B83:AC83").Value2
$ws.Range("B82:AC82").Value2 = $rowB
$ws.Range("B83:AC83").Value2 = $rowA

$rowA = $ws.Range("B104:AC104").Value2
$rowB = $ws.Range("B105:AC105").Value2
$ws.Range("B104:AC104").Value2 = $rowB
$ws.Range("B105:AC105").Value2 = $rowA

$rowA = $ws.Range("B129:AC129").Value2
$rowB = $ws.Range("B130:AC130").Value2
$ws.Range("B129:AC129").Value2 = $rowB
$ws.Range("B130:AC130").Value2 = $rowA

$rowA = $ws.Range("B142:AC142").Value2
$rowB = $ws.Range("B143:AC143").Value2
$ws.Range("B142:AC142").Value2 = $rowB
$ws.Range("B143:AC143").Value2 = $rowA

$rowA = $ws.Range("B164:AC164").Value2
$rowB = $ws.Range("B165:AC165").Value2
$ws.Range("B164:AC164").Value2 = $rowB
$ws.Range("B165:AC165").Value2 = $rowA

$rowA = $ws.Range("B167:AC167").Value2
$rowB = $ws.Range("B168:AC168").Value2
$ws.Range("B167:AC167").Value2 = $rowB
$ws.Range("B168:AC168").Value2 = $rowA

$rowA = $ws.Range("B171:AC171").Value2
$rowB = $ws.Range("B172:AC172").Value2
$ws.Range("B171:AC171").Value2 = $rowB
$ws.Range("B172:AC172").Value2 = $rowA

$rowA = $ws.Range("B177:AC177").Value2
$rowB = $ws.Range("B178:AC178").Value2
$ws.Range("B177:AC177").Value2 = $rowB
$ws.Range("B178:AC178").Value2 = $rowA

# Rotate rows 180 -> 181 -> 183 -> 180 (3-cycle).
$row180 = $ws.Range("B180:AC180").Value2
$row181 = $ws.Range("B181:AC181").Value2
$row183 = $ws.Range("B183:AC183").Value2
$ws.Range("B181:AC181").Value2 = $row180
$ws.Range("B183:AC183").Value2 = $row181
$ws.Range("B180:AC180").Value2 = $row183
